# Generate Report for Handback
# ---------------------------------------------------------------------------
# This applies the "handback" update to the localization-status workbook:
#  - Overview sheet: zh-cn / de-de status cells move from "Ready for handoff"
#    to "Handed back: in sync with en-US" (and widen those columns).
#  - zh-cn / de-de detail sheets: the "Latest Target File" column becomes a
#    hyperlink to the source markdown file, "Latest Handback File" is filled
#    in with the generated xlf file name, and "Latest Handback DateTime" is
#    stamped with the handback timestamp. Related columns are widened to fit.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/9bbfdef86382df90f0c48f5f9b3627ba4eafbd27/e2e/d6a5534b-124a-4535-8f3c-cfb62ed6460e.md"
$mdName = "d6a5534b-124a-4535-8f3c-cfb62ed6460e.md"
$statusText = "Handed back: in sync with en-US"

# Matches the underline + font colour used by the workbook's existing
# hyperlink style (single underline, RGB FF6495ED == BGR 15570276).
$hyperlinkUnderline = 2
$hyperlinkColor = 15570276

# -----------------------------------------------------------------------
# Overview sheet
# -----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# -----------------------------------------------------------------------
# zh-cn sheet
# -----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17
$wsZhCn.Columns.Item(11).ColumnWidth = 39.17

$wsZhCn.Range("C2").Value = $statusText

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J2"), $mdUrl, "", "", $mdName)
$wsZhCn.Range("J2").Font.Underline = $hyperlinkUnderline
$wsZhCn.Range("J2").Font.Color = $hyperlinkColor

$wsZhCn.Range("K2").Value = "d6a5534b-124a-4535-8f3c-cfb62ed6460e.a99fa954c34e681f7b56e8fb81b2a7a3bbc97707.zh-cn.xlf"
$wsZhCn.Range("L2").Value = "2017-01-03 05:26:10"

# -----------------------------------------------------------------------
# de-de sheet
# -----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17
$wsDeDe.Columns.Item(11).ColumnWidth = 39.17

$wsDeDe.Range("C2").Value = $statusText

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J2"), $mdUrl, "", "", $mdName)
$wsDeDe.Range("J2").Font.Underline = $hyperlinkUnderline
$wsDeDe.Range("J2").Font.Color = $hyperlinkColor

$wsDeDe.Range("K2").Value = "d6a5534b-124a-4535-8f3c-cfb62ed6460e.a99fa954c34e681f7b56e8fb81b2a7a3bbc97707.de-de.xlf"
$wsDeDe.Range("L2").Value = "2017-01-03 05:26:21"

Write-Host "Handback report generated"
